$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "24÷7=3, 3") { Write-Host "MISMATCH row=1 col=1 got=" $r.Text }
$r.Text = "27÷4=6, 3"

$cell = $t.Cell(1,2)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "27÷8=3, 3") { Write-Host "MISMATCH row=1 col=2 got=" $r.Text }
$r.Text = "53÷4=13, 1"

$cell = $t.Cell(1,3)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "53÷2=26, 1") { Write-Host "MISMATCH row=1 col=3 got=" $r.Text }
$r.Text = "33÷2=16, 1"

$cell = $t.Cell(1,4)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "60÷6=10, 0") { Write-Host "MISMATCH row=1 col=4 got=" $r.Text }
$r.Text = "31÷9=3, 4"

$cell = $t.Cell(1,5)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "81÷9=9, 0") { Write-Host "MISMATCH row=1 col=5 got=" $r.Text }
$r.Text = "45÷6=7, 3"

$cell = $t.Cell(5,1)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "92÷3=30, 2") { Write-Host "MISMATCH row=5 col=1 got=" $r.Text }
$r.Text = "79÷4=19, 3"

$cell = $t.Cell(5,2)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "93÷9=10, 3") { Write-Host "MISMATCH row=5 col=2 got=" $r.Text }
$r.Text = "73÷7=10, 3"

$cell = $t.Cell(5,3)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "60÷8=7, 4") { Write-Host "MISMATCH row=5 col=3 got=" $r.Text }
$r.Text = "42÷7=6, 0"

$cell = $t.Cell(5,4)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "15÷9=1, 6") { Write-Host "MISMATCH row=5 col=4 got=" $r.Text }
$r.Text = "31÷4=7, 3"

$cell = $t.Cell(5,5)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "71÷5=14, 1") { Write-Host "MISMATCH row=5 col=5 got=" $r.Text }
$r.Text = "64÷8=8, 0"

$cell = $t.Cell(9,1)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "54÷2=27, 0") { Write-Host "MISMATCH row=9 col=1 got=" $r.Text }
$r.Text = "68÷3=22, 2"

$cell = $t.Cell(9,2)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "34÷6=5, 4") { Write-Host "MISMATCH row=9 col=2 got=" $r.Text }
$r.Text = "57÷2=28, 1"

$cell = $t.Cell(9,3)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "22÷9=2, 4") { Write-Host "MISMATCH row=9 col=3 got=" $r.Text }
$r.Text = "47÷8=5, 7"

$cell = $t.Cell(9,4)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "70÷3=23, 1") { Write-Host "MISMATCH row=9 col=4 got=" $r.Text }
$r.Text = "46÷7=6, 4"

$cell = $t.Cell(9,5)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "38÷4=9, 2") { Write-Host "MISMATCH row=9 col=5 got=" $r.Text }
$r.Text = "72÷7=10, 2"

$cell = $t.Cell(13,1)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "77÷7=11, 0") { Write-Host "MISMATCH row=13 col=1 got=" $r.Text }
$r.Text = "29÷9=3, 2"

$cell = $t.Cell(13,2)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "32÷4=8, 0") { Write-Host "MISMATCH row=13 col=2 got=" $r.Text }
$r.Text = "34÷5=6, 4"

$cell = $t.Cell(13,3)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "10÷9=1, 1") { Write-Host "MISMATCH row=13 col=3 got=" $r.Text }
$r.Text = "80÷4=20, 0"

$cell = $t.Cell(13,4)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "60÷6=10, 0") { Write-Host "MISMATCH row=13 col=4 got=" $r.Text }
$r.Text = "74÷7=10, 4"

$cell = $t.Cell(13,5)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "74÷4=18, 2") { Write-Host "MISMATCH row=13 col=5 got=" $r.Text }
$r.Text = "75÷3=25, 0"

$cell = $t.Cell(17,1)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "36÷9=4, 0") { Write-Host "MISMATCH row=17 col=1 got=" $r.Text }
$r.Text = "32÷4=8, 0"

$cell = $t.Cell(17,2)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "73÷9=8, 1") { Write-Host "MISMATCH row=17 col=2 got=" $r.Text }
$r.Text = "35÷5=7, 0"

$cell = $t.Cell(17,3)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "31÷2=15, 1") { Write-Host "MISMATCH row=17 col=3 got=" $r.Text }
$r.Text = "83÷3=27, 2"

$cell = $t.Cell(17,4)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "11÷2=5, 1") { Write-Host "MISMATCH row=17 col=4 got=" $r.Text }
$r.Text = "19÷7=2, 5"

$cell = $t.Cell(17,5)
$full = $cell.Range
$r = $d.Range($full.Start, $full.End - 1)
if ($r.Text -ne "72÷5=14, 2") { Write-Host "MISMATCH row=17 col=5 got=" $r.Text }
$r.Text = "37÷4=9, 1"
